$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'45.261.78"
$ws.Range("E2").Value = "  +3.68%  "
$ws.Range("D3").Value = "'2.428.12"
$ws.Range("E3").Value = "  +0.39%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").Value = "'317.44"
$ws.Range("E5").Value = "  +3.45%  "
$ws.Range("D6").Value = "'102.65"
$ws.Range("E6").Value = "  +5.19%  "
$ws.Range("E7").Value = "  +1.57%  "
$ws.Range("E8").Value = "  -0.13%  "
$ws.Range("E9").Value = "  +7.13%  "
$ws.Range("D10").Value = "'35.61"
$ws.Range("E10").Value = "  +1.60%  "
$ws.Range("E11").Value = "  +0.70%  "
$ws.Range("E12").Value = "  -2.20%  "
$ws.Range("D13").Value = "'18.14"
$ws.Range("D14").Value = "'7.05"
$ws.Range("E14").Value = "  +2.09%  "
$ws.Range("D15").Value = "'2.809.11"
$ws.Range("E15").Value = "  +0.89%  "
$ws.Range("D16").Value = "'2.433.89"
$ws.Range("E16").Value = "  +1.96%  "
$ws.Range("D17").Value = "'0.843"
$ws.Range("E17").Value = "  +2.01%  "
$ws.Range("D18").Value = "'45.160.96"
$ws.Range("E18").Value = "  +3.50%  "
$ws.Range("D19").Value = "'12.23"
$ws.Range("E19").Value = "  +0.83%  "
$ws.Range("E20").Value = "  -0.88%  "
$ws.Range("E21").Value = "  +2.17%  "
$ws.Range("E22").Value = "  +0.75%  "
$ws.Range("D23").Value = "'243.98"
$ws.Range("E23").Value = "  +2.21%  "
$ws.Range("E24").Value = "  +0.44%  "
$ws.Range("E25").Value = "  +1.78%  "
$ws.Range("E26").Value = "  -0.03%  "
$ws.Range("D27").Value = "'25.57"
$ws.Range("E27").Value = "  +2.30%  "
$ws.Range("D28").Value = "'9.59"
$ws.Range("E28").Value = "  +1.56%  "
$ws.Range("E29").Value = "  -11.96%  "
$ws.Range("B30").Value = "InjectiveProtocol"
$ws.Range("C30").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D30").Value = "'32.94"
$ws.Range("E30").Value = "  +1.67%  "
$ws.Range("B31").Value = "OKB"
$ws.Range("C31").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D31").Value = "'49.07"
$ws.Range("E31").Value = "  +1.99%  "
$ws.Range("D32").Value = "'20.34"
$ws.Range("E32").Value = "  +10.54%  "
$ws.Range("E33").Value = "  +4.49%  "
$ws.Range("E34").Value = "  +1.26%  "
$ws.Range("E35").Value = "  +0.28%  "
$ws.Range("D36").Value = "'0.0765"
$ws.Range("E36").Value = "  +1.58%  "
$ws.Range("E37").Value = "  -1.15%  "
$ws.Range("D38").Value = "'4.45"
$ws.Range("E38").Value = "  +1.04%  "
$ws.Range("D39").Value = "'2.85"
$ws.Range("E39").Value = "  -2.09%  "
$ws.Range("D40").Value = "'125.58"
$ws.Range("E40").Value = "  -3.64%  "
$ws.Range("D41").Value = "'2.22"
$ws.Range("E41").Value = "  -2.75%  "
$ws.Range("E42").Value = "  +0.59%  "
$ws.Range("D43").Value = "'20.57"
$ws.Range("E43").Value = "  -2.83%  "
$ws.Range("E44").Value = "  +1.99%  "
$ws.Range("D45").Value = "'1.937.04"
$ws.Range("E45").Value = "  -0.59%  "
$ws.Range("E46").Value = "  -2.87%  "
$ws.Range("D47").Value = "'2.92"
$ws.Range("E47").Value = "  +3.14%  "
$ws.Range("E48").Value = "  +16.06%  "
$ws.Range("E49").Value = "  -1.40%  "
$ws.Range("D50").Value = "'76.47"
$ws.Range("E50").Value = "  +5.62%  "
$ws.Range("D51").Value = "'53.86"
$ws.Range("E51").Value = "  +1.88%  "
